{"js": "// The diff moves \"2022 \" from after \"\u03a4\u03b1\u03cd\u03c1\u03bf\u03c5\" to the start of the sentence,\n// in all four occurrences of the sentence across the document body.\nconst oldText =\n  \"\u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03a4\u03b1\u03cd\u03c1\u03bf\u03c5 2022: 16-25 \u0399\u03b1\u03bd\u03bf\u03c5\u03b1\u03c1\u03af\u03bf\u03c5\";\nconst newText =\n  \"2022 \u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03a4\u03b1\u03cd\u03c1\u03bf\u03c5: 16-25 \u0399\u03b1\u03bd\u03bf\u03c5\u03b1\u03c1\u03af\u03bf\u03c5\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (const range of results.items) {\n  range.insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# The diff moves \"2022 \" from after \"\u03a4\u03b1\u03cd\u03c1\u03bf\u03c5\" to the start of the sentence,\n# in all four occurrences of the sentence across the document body.\n$d = $word.ActiveDocument\n\n$oldText = \"\u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03a4\u03b1\u03cd\u03c1\u03bf\u03c5 2022: 16-25 \u0399\u03b1\u03bd\u03bf\u03c5\u03b1\u03c1\u03af\u03bf\u03c5\"\n$newText = \"2022 \u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03a4\u03b1\u03cd\u03c1\u03bf\u03c5: 16-25 \u0399\u03b1\u03bd\u03bf\u03c5\u03b1\u03c1\u03af\u03bf\u03c5\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n"}
